# Update cryptocurrency price/volume figures in the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text happens to look like a plain number ---
# (e.g. "300.77", "0.0299") need special handling so Excel keeps
# storing them as literal text instead of silently converting them
# to a numeric value (which would also reformat/round them).
# We temporarily force a "Text" number format, assign the value,
# then restore each cell's original style so no visible
# formatting/style change is left behind.

$numericLookingCells = @(
    "D5"
    "D6"
    "D7"
    "D9"
    "D10"
    "D11"
    "D12"
    "D16"
    "D17"
    "D19"
    "D21"
    "D22"
    "D23"
    "D24"
    "D25"
    "D26"
    "D27"
    "D29"
    "D30"
    "D31"
    "D33"
    "D34"
    "D35"
    "D37"
    "D39"
    "D40"
    "D41"
    "D42"
    "D45"
    "D49"
    "D51"
)

$savedStyles = @{}
foreach ($addr in $numericLookingCells) {
    $savedStyles[$addr] = $ws.Range($addr).Style
}
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$numericLookingValues = @{
    "D5" = '300.77'
    "D6" = '98.44'
    "D7" = '0.564'
    "D9" = '0.508'
    "D10" = '34.49'
    "D11" = '0.0788'
    "D12" = '7.12'
    "D16" = '0.824'
    "D17" = '13.73'
    "D19" = '12.73'
    "D21" = '6.05'
    "D22" = '66.77'
    "D23" = '243.13'
    "D24" = '2.78'
    "D25" = '0.999'
    "D26" = '1.93'
    "D27" = '39.72'
    "D29" = '9.74'
    "D30" = '3.77'
    "D31" = '20.91'
    "D33" = '5.54'
    "D34" = '146.84'
    "D35" = '0.0773'
    "D37" = '1.92'
    "D39" = '15.09'
    "D40" = '3.88'
    "D41" = '0.0299'
    "D42" = '3.21'
    "D45" = '91.98'
    "D49" = '98.94'
    "D51" = '68.78'
}
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Value = $numericLookingValues[$addr]
}
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = $savedStyles[$addr]
}

# --- Remaining Price cells (values that already are not plain numbers,
# e.g. "45.982.00" which Excel can only interpret as text) ---
$ws.Range("D2").Value = '45.982.00'
$ws.Range("D3").Value = '2.377.40'
$ws.Range("D14").Value = '2.740.36'
$ws.Range("D15").Value = '2.384.85'
$ws.Range("D18").Value = '45.886.61'
$ws.Range("D43").Value = '1.933.34'
$ws.Range("D50").Value = '2.611.42'

# --- Volume(1h) percentage cells (always text, never numeric-looking) ---
$ws.Range("E2").Value = '  -1.69%  '
$ws.Range("E3").Value = '  +3.21%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("E6").Value = '  -3.19%  '
$ws.Range("E7").Value = '  -1.12%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  -4.74%  '
$ws.Range("E10").Value = '  -6.15%  '
$ws.Range("E11").Value = '  -2.02%  '
$ws.Range("E12").Value = '  -3.40%  '
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("E14").Value = '  +3.18%  '
$ws.Range("E15").Value = '  +3.59%  '
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("E17").Value = '  -2.06%  '
$ws.Range("E18").Value = '  -1.77%  '
$ws.Range("E19").Value = '  -5.35%  '
$ws.Range("E20").Value = '  +0.58%  '
$ws.Range("E21").Value = '  -1.24%  '
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("E23").Value = '  -2.22%  '
$ws.Range("E24").Value = '  -5.92%  '
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("E26").Value = '  -2.36%  '
$ws.Range("E27").Value = '  -11.32%  '
$ws.Range("E28").Value = '  -3.44%  '
$ws.Range("E29").Value = '  -1.63%  '
$ws.Range("E30").Value = '  +20.24%  '
$ws.Range("E31").Value = '  +3.65%  '
$ws.Range("E32").Value = '  +6.82%  '
$ws.Range("E33").Value = '  -4.57%  '
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("E35").Value = '  -3.37%  '
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("E37").Value = '  +5.90%  '
$ws.Range("E38").Value = '  -2.66%  '
$ws.Range("E39").Value = '  -4.96%  '
$ws.Range("E40").Value = '  -4.30%  '
$ws.Range("E41").Value = '  -1.99%  '
$ws.Range("E42").Value = '  -9.09%  '
$ws.Range("E43").Value = '  +3.58%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("E45").Value = '  +4.58%  '
$ws.Range("E46").Value = '  -9.35%  '
$ws.Range("E47").Value = '  +4.91%  '
$ws.Range("E48").Value = '  -5.16%  '
$ws.Range("E49").Value = '  +1.51%  '
$ws.Range("E50").Value = '  +3.15%  '
$ws.Range("E51").Value = '  -7.70%  '

